$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2538
$ws.Range("E2").Value = 171
$ws.Range("F2").Value = 171
$ws.Range("G2").Value = 181
$ws.Range("H2").Value = 138
$ws.Range("I2").Value = 138
$ws.Range("K2").Value = 1396
$ws.Range("L2").Value = 790
$ws.Range("M2").Value = 606
$ws.Range("N2").Value = 606
$ws.Range("P2").Value = 22
$ws.Range("Q2").Value = -30
$ws.Range("R2").Value = -154
$ws.Range("S2").Value = 56
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = -30
$ws.Range("V2").Value = 650
$ws.Range("W2").Value = 6.75
$ws.Range("X2").Value = 5.42
$ws.Range("Y2").Value = 25.63
$ws.Range("Z2").Value = 11.24
$ws.Range("AA2").Value = 130.39
$ws.Range("AB2").Value = 2719.12
$ws.Range("AC2").Value = 2307
$ws.Range("AE2").Value = 10155
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 5168480
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AH2").ClearContents()

# Row 3
$ws.Range("D3").Value = 2969
$ws.Range("E3").Value = 250
$ws.Range("F3").Value = 250
$ws.Range("G3").Value = 343
$ws.Range("H3").Value = 258
$ws.Range("I3").Value = 258
$ws.Range("K3").Value = 2283
$ws.Range("L3").Value = 1692
$ws.Range("M3").Value = 591
$ws.Range("N3").Value = 588
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 34
$ws.Range("Q3").Value = 149
$ws.Range("R3").Value = -196
$ws.Range("S3").Value = 9
$ws.Range("T3").Value = 10
$ws.Range("U3").Value = 139
$ws.Range("V3").Value = 1292
$ws.Range("W3").Value = 8.41
$ws.Range("X3").Value = 8.699999999999999
$ws.Range("Y3").Value = 43.26
$ws.Range("Z3").Value = 14.04
$ws.Range("AA3").Value = 286.37
$ws.Range("AB3").Value = 2414.48
$ws.Range("AC3").Value = 4324
$ws.Range("AE3").Value = 9478
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 5405345
$ws.Range("J3").ClearContents()
$ws.Range("AD3").ClearContents()
$ws.Range("AH3").ClearContents()

# Row 4
$ws.Range("D4").Value = 3170
$ws.Range("E4").Value = 280
$ws.Range("F4").Value = 280
$ws.Range("G4").Value = 250
$ws.Range("H4").Value = 161
$ws.Range("I4").Value = 159
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 2427
$ws.Range("L4").Value = 1531
$ws.Range("M4").Value = 896
$ws.Range("N4").Value = 891
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = 41
$ws.Range("Q4").Value = 134
$ws.Range("R4").Value = 9
$ws.Range("S4").Value = -102
$ws.Range("T4").Value = 114
$ws.Range("U4").Value = 20
$ws.Range("V4").Value = 1162
$ws.Range("W4").Value = 8.83
$ws.Range("X4").Value = 5.07
$ws.Range("Y4").Value = 21.53
$ws.Range("Z4").Value = 6.83
$ws.Range("AA4").Value = 170.85
$ws.Range("AB4").Value = 2657.7
$ws.Range("AC4").Value = 2455
$ws.Range("AE4").Value = 12961
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 6871990
$ws.Range("AD4").ClearContents()
$ws.Range("AH4").ClearContents()

# Row 5
$ws.Range("D5").Value = 3284
$ws.Range("E5").Value = 151
$ws.Range("F5").Value = 151
$ws.Range("G5").Value = 69
$ws.Range("H5").Value = 44
$ws.Range("I5").Value = 45
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 2853
$ws.Range("L5").Value = 1637
$ws.Range("M5").Value = 1215
$ws.Range("N5").Value = 1211
$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 50
$ws.Range("Q5").Value = -180
$ws.Range("R5").Value = -250
$ws.Range("S5").Value = 486
$ws.Range("T5").Value = 153
$ws.Range("U5").Value = -332
$ws.Range("V5").Value = 1245
$ws.Range("W5").Value = 4.59
$ws.Range("X5").Value = 1.35
$ws.Range("Y5").Value = 4.29
$ws.Range("Z5").Value = 1.68
$ws.Range("AA5").Value = 134.76
$ws.Range("AB5").Value = 3013.19
$ws.Range("AC5").Value = 547
$ws.Range("AD5").Value = 23.93
$ws.Range("AE5").Value = 14871
$ws.Range("AF5").Value = 0.88
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 1.15
$ws.Range("AI5").Value = 27.11
$ws.Range("AJ5").Value = 8000000

# Row 6
$ws.Range("D6").Value = 3133
$ws.Range("E6").Value = 108
$ws.Range("F6").Value = 108
$ws.Range("G6").Value = 58
$ws.Range("H6").Value = 46
$ws.Range("I6").Value = 48
$ws.Range("K6").Value = 2875
$ws.Range("L6").Value = 1610
$ws.Range("M6").Value = 1265
$ws.Range("N6").Value = 1260
$ws.Range("P6").Value = 50
$ws.Range("Q6").Value = 56
$ws.Range("R6").Value = -55
$ws.Range("S6").Value = -84
$ws.Range("T6").Value = 119
$ws.Range("U6").Value = -63
$ws.Range("V6").Value = 1250
$ws.Range("W6").Value = 3.43
$ws.Range("X6").Value = 1.47
$ws.Range("Y6").Value = 3.88
$ws.Range("Z6").Value = 1.61
$ws.Range("AA6").Value = 127.28
$ws.Range("AB6").Value = 3100.37
$ws.Range("AC6").Value = 555
$ws.Range("AD6").Value = 16.67
$ws.Range("AE6").Value = 15468
$ws.Range("AF6").Value = 0.6
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 2.7
$ws.Range("AI6").Value = 42.46
$ws.Range("AJ6").Value = 8000000

# Row 7
$ws.Range("D7").Value = 3144
$ws.Range("E7").Value = 203
$ws.Range("G7").Value = 133
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 99
$ws.Range("K7").Value = 2915
$ws.Range("L7").Value = 1553
$ws.Range("M7").Value = 1361
$ws.Range("N7").Value = 1356
$ws.Range("P7").Value = 50
$ws.Range("Q7").Value = 137
$ws.Range("R7").Value = -115
$ws.Range("S7").Value = -38
$ws.Range("T7").Value = 0
$ws.Range("W7").Value = 6.46
$ws.Range("X7").Value = 3.18
$ws.Range("Y7").Value = 7.57
$ws.Range("Z7").Value = 3.45
$ws.Range("AA7").Value = 114.11
$ws.Range("AC7").Value = 1146
$ws.Range("AD7").Value = 12.57
$ws.Range("AE7").Value = 16649
$ws.Range("AF7").Value = 0.86
$ws.Range("AG7").Value = 250
$ws.Range("AH7").Value = 1.74
$ws.Range("AI7").Value = 20.2
$ws.Range("U7").ClearContents()

# Row 8
$ws.Range("D8").Value = 3294
$ws.Range("E8").Value = 216
$ws.Range("G8").Value = 136
$ws.Range("H8").Value = 103
$ws.Range("I8").Value = 103
$ws.Range("K8").Value = 3098
$ws.Range("L8").Value = 1637
$ws.Range("M8").Value = 1461
$ws.Range("N8").Value = 1455
$ws.Range("P8").Value = 50
$ws.Range("Q8").Value = 387
$ws.Range("R8").Value = -120
$ws.Range("S8").Value = 63
$ws.Range("T8").Value = 0
$ws.Range("W8").Value = 6.56
$ws.Range("X8").Value = 3.13
$ws.Range("Y8").Value = 7.33
$ws.Range("Z8").Value = 3.43
$ws.Range("AA8").Value = 112.05
$ws.Range("AC8").Value = 1192
$ws.Range("AD8").Value = 12.08
$ws.Range("AE8").Value = 17865
$ws.Range("AF8").Value = 0.8100000000000001
$ws.Range("AG8").Value = 250
$ws.Range("AH8").Value = 1.74
$ws.Range("AI8").Value = 19.42
$ws.Range("U8").ClearContents()

# Row 9
$ws.Range("D9").Value = 3545
$ws.Range("E9").Value = 268
$ws.Range("G9").Value = 188
$ws.Range("H9").Value = 143
$ws.Range("I9").Value = 142
$ws.Range("K9").Value = 3232
$ws.Range("L9").Value = 1632
$ws.Range("M9").Value = 1600
$ws.Range("N9").Value = 1593
$ws.Range("P9").Value = 50
$ws.Range("Q9").Value = 119
$ws.Range("R9").Value = -124
$ws.Range("S9").Value = -37
$ws.Range("T9").Value = 0
$ws.Range("W9").Value = 7.56
$ws.Range("X9").Value = 4.03
$ws.Range("Y9").Value = 9.32
$ws.Range("Z9").Value = 4.52
$ws.Range("AA9").Value = 102
$ws.Range("AC9").Value = 1643
$ws.Range("AD9").Value = 8.76
$ws.Range("AE9").Value = 19559
$ws.Range("AF9").Value = 0.74
$ws.Range("AG9").Value = 250
$ws.Range("AH9").Value = 1.74
$ws.Range("AI9").Value = 14.09
$ws.Range("U9").ClearContents()
